# Update generated counts/prices in the "展览" (exhibition) and
# "全部类型" (all types) sheets to reflect the latest scrape output.
# Corresponds to commit: "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Sheet "展览": cell -> new value
$wsExhibition.Range("F2").Value  = 6596
$wsExhibition.Range("F4").Value  = 412
$wsExhibition.Range("F9").Value  = 90
$wsExhibition.Range("F11").Value = 166
$wsExhibition.Range("F12").Value = 389
$wsExhibition.Range("F13").Value = 1273
$wsExhibition.Range("F15").Value = 3277
$wsExhibition.Range("G16").Value = 45
$wsExhibition.Range("F17").Value = 209
$wsExhibition.Range("F18").Value = 1912
$wsExhibition.Range("F19").Value = 25
$wsExhibition.Range("F21").Value = 121

# Sheet "全部类型": same logical rows, shifted by one (extra row present)
$wsAllTypes.Range("F2").Value  = 6596
$wsAllTypes.Range("F4").Value  = 412
$wsAllTypes.Range("F10").Value = 90
$wsAllTypes.Range("F12").Value = 166
$wsAllTypes.Range("F13").Value = 389
$wsAllTypes.Range("F14").Value = 1273
$wsAllTypes.Range("F16").Value = 3277
$wsAllTypes.Range("G17").Value = 45
$wsAllTypes.Range("F18").Value = 209
$wsAllTypes.Range("F19").Value = 1912
$wsAllTypes.Range("F20").Value = 25
$wsAllTypes.Range("F22").Value = 121
